$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.434.29'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '3.072.19'
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = '234.59'
$ws.Range('E5').Value = '  +8.57%  '
$ws.Range('D6').Value = '617.68'
$ws.Range('E6').Value = '  -1.36%  '
$ws.Range('E7').Value = '  -8.92%  '
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Value = '3.069.64'
$ws.Range('E10').Value = '  -2.63%  '
$ws.Range('D11').Value = '0.712'
$ws.Range('E11').Value = '  -6.88%  '
$ws.Range('D12').Value = '0.197'
$ws.Range('E12').Value = '  -1.88%  '
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').Value = '35.13'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D15').Value = '89.546.62'
$ws.Range('E15').Value = '  -1.26%  '
$ws.Range('D16').Value = '5.35'
$ws.Range('E16').Value = '  -6.52%  '
$ws.Range('D17').Value = '3.641.50'
$ws.Range('E17').Value = '  -2.63%  '
$ws.Range('D18').Value = '3.100.03'
$ws.Range('E18').Value = '  -3.38%  '
$ws.Range('D19').Value = '3.79'
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').Value = '0.0000211'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').Value = '13.70'
$ws.Range('E21').Value = '  -6.36%  '
$ws.Range('D22').Value = '431.44'
$ws.Range('E22').Value = '  -9.31%  '
$ws.Range('D23').Value = '5.38'
$ws.Range('E23').Value = '  +3.86%  '
$ws.Range('D24').Value = '8.73'
$ws.Range('E24').Value = '  -4.67%  '
$ws.Range('D25').Value = '5.57'
$ws.Range('E25').Value = '  -6.73%  '
$ws.Range('D26').Value = '86.82'
$ws.Range('E26').Value = '  -8.29%  '
$ws.Range('D27').Value = '11.70'
$ws.Range('E27').Value = '  -5.52%  '
$ws.Range('E28').Value = '  -2.21%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '9.02'
$ws.Range('E31').Value = '  -3.38%  '
$ws.Range('D32').Value = '0.155'
$ws.Range('E32').Value = '  -4.72%  '
$ws.Range('E33').Value = '  -8.89%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '25.45'
$ws.Range('E34').Value = '  -7.47%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.150'
$ws.Range('E35').Value = '  +2.70%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').Value = '7.09'
$ws.Range('E36').Value = '  +1.94%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').Value = '3.67'
$ws.Range('E37').Value = '  +1.74%  '
$ws.Range('D38').Value = '494.22'
$ws.Range('E38').Value = '  -4.80%  '
$ws.Range('D39').Value = '1.87'
$ws.Range('E39').Value = '  -3.61%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.0906'
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = '1.25'
$ws.Range('E41').Value = '  -4.85%  '
$ws.Range('E42').Value = '  +54.42%  '
$ws.Range('D43').Value = '22.07'
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('E45').Value = '  -7.34%  '
$ws.Range('D46').Value = '151.59'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('E47').Value = '  -7.08%  '
$ws.Range('D48').Value = '0.673'
$ws.Range('E48').Value = '  -8.72%  '
$ws.Range('D49').Value = '44.33'
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('D50').Value = '0.997'
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('E51').Value = '  -5.47%  '
